# VAS , KEY account changes
#  - Rename the acBandsPalletsSum / dBandCasesSum headers on lessThan100Cases
#  - Update the I/J (volume-band) figures for a number of rows
#  - Add a new "volBands" worksheet summarising counts/skus/shipments per band

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the two volume-band summary headers -----------------------
$ws.Range("I1").Value = "abBandsPalletsSum"
$ws.Range("J1").Value = "cdBandCasesSum"

# --- Update the per-row volume-band figures (columns I & J) -----------
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 30

$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 27

$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 77

$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 27

$ws.Range("I18").Value = 2
$ws.Range("J18").Value = 55

$ws.Range("I24").Value = 8
$ws.Range("J24").Value = 42

$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 31

$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 167

$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 21

$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 113

$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 37

$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 20

$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 82

$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 15

# --- Add the new volBands worksheet, right after lessThan100Cases -----
$volWs = $wb.Worksheets.Add($null, $ws)
$volWs.Name = "volBands"

# Header row
$volWs.Range("A1").Value = "volBand"
$volWs.Range("B1").Value = "count"
$volWs.Range("C1").Value = "unitsSum"
$volWs.Range("D1").Value = "casesSum"
$volWs.Range("E1").Value = "palletsSum"
$volWs.Range("F1").Value = "pallEquivSum"
$volWs.Range("G1").Value = "skus"
$volWs.Range("H1").Value = "shipments"

# Band A
$volWs.Range("A2").Value = "A"
$volWs.Range("B2").Value = 1362
$volWs.Range("C2").Value = 1335269
$volWs.Range("D2").Value = 215304
$volWs.Range("E2").Value = 9665
$volWs.Range("F2").Value = 8972.03
$volWs.Range("G2").Value = 921
$volWs.Range("H2").Value = 303

# Band B
$volWs.Range("A3").Value = "B"
$volWs.Range("B3").Value = 2657
$volWs.Range("C3").Value = 721043
$volWs.Range("D3").Value = 105543
$volWs.Range("E3").Value = 5944
$volWs.Range("F3").Value = 4399.5
$volWs.Range("G3").Value = 2119
$volWs.Range("H3").Value = 386

# Band C
$volWs.Range("A4").Value = "C"
$volWs.Range("B4").Value = 3116
$volWs.Range("C4").Value = 376102
$volWs.Range("D4").Value = 52903
$volWs.Range("E4").Value = 3116
$volWs.Range("F4").Value = 2206.11
$volWs.Range("G4").Value = 2719
$volWs.Range("H4").Value = 374

# Band D
$volWs.Range("A5").Value = "D"
$volWs.Range("B5").Value = 8445
$volWs.Range("C5").Value = 283116
$volWs.Range("D5").Value = 39133
$volWs.Range("E5").Value = 8444
$volWs.Range("F5").Value = 1634.35
$volWs.Range("G5").Value = 6871
$volWs.Range("H5").Value = 470

$ws.Activate()
